$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1106.1
$ws.Range("I62").Value = 953.75
$ws.Range("J62").Value = 1207.6666
$ws.Range("K62").Value = 953.75
$ws.Range("L62").Value = 1207.6666
$ws.Range("M62").Value = -329.75
$ws.Range("N62").Value = -2455.6666
$ws.Range("H65").Value = 1106.1
$ws.Range("I65").Value = 953.75
$ws.Range("J65").Value = 1207.6666
$ws.Range("K65").Value = 4768.75
$ws.Range("L65").Value = 6038.333000000001
$ws.Range("M65").Value = -1648.75
$ws.Range("N65").Value = -12278.333
$ws.Range("H100").Value = 2006
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2006
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = $null
$ws.Range("M100").Value = 2006
$ws.Range("N100").Value = -3088
$ws.Range("H106").Value = 52749.5
$ws.Range("I106").Value = 58110.555
$ws.Range("K106").Value = 58110.555
$ws.Range("M106").Value = -57479.555
$ws.Range("H116").Value = 3330.5652
$ws.Range("I116").Value = 2449.7856
$ws.Range("J116").Value = 4700.6665
$ws.Range("K116").Value = 2449.7856
$ws.Range("L116").Value = 4700.6665
$ws.Range("M116").Value = 992.2143999999998
$ws.Range("N116").Value = -11584.6665
$ws.Range("H130").Value = 44000
$ws.Range("J130").Value = 44000
$ws.Range("L130").Value = 44000
$ws.Range("N130").Value = -54040
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = $null
$ws.Range("N134").Value = 0
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8751.8125
$ws.Range("I32").Value = 9272.983
$ws.Range("J32").Value = 7078.579
$ws.Range("K32").Value = 9272.983
$ws.Range("L32").Value = 7078.579
$ws.Range("M32").Value = -8985.983
$ws.Range("N32").Value = -7652.579
$ws.Range("H102").Value = 2473.6365
$ws.Range("I102").Value = 2245.5557
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2245.5557
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -623.5556999999999
$ws.Range("N102").Value = -6744
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 34503.6
$ws.Range("J76").Value = 34503.6
$ws.Range("L76").Value = 34503.6
$ws.Range("N76").Value = -35133.6
$ws.Range("H79").Value = 34503.6
$ws.Range("J79").Value = 34503.6
$ws.Range("L79").Value = 34503.6
$ws.Range("N79").Value = -36687.6
$ws.Range("H132").Value = 58393.168
$ws.Range("J132").Value = 58393.168
$ws.Range("L132").Value = 58393.168
$ws.Range("N132").Value = -68513.16800000001
$ws.Range("H134").Value = 2081.5
$ws.Range("I134").Value = 1367.5
$ws.Range("J134").Value = 3628.5
$ws.Range("K134").Value = 4102.5
$ws.Range("L134").Value = 10885.5
$ws.Range("M134").Value = -1567.5
$ws.Range("N134").Value = -15955.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4330616.5
$ws.Range("I31").Value = 1345.2667
$ws.Range("J31").Value = 19610398
$ws.Range("K31").Value = 1345.2667
$ws.Range("L31").Value = 19610398
$ws.Range("M31").Value = -1050.2667
$ws.Range("N31").Value = -19610988
$ws.Range("H34").Value = 4330616.5
$ws.Range("I34").Value = 1345.2667
$ws.Range("J34").Value = 19610398
$ws.Range("K34").Value = 1345.2667
$ws.Range("L34").Value = 19610398
$ws.Range("M34").Value = -1143.2667
$ws.Range("N34").Value = -19610802
$ws.Range("H99").Value = 1922.6154
$ws.Range("I99").Value = 1527.7142
$ws.Range("J99").Value = 2383.3333
$ws.Range("K99").Value = 1527.7142
$ws.Range("L99").Value = 2383.3333
$ws.Range("M99").Value = -29.71419999999989
$ws.Range("N99").Value = -5379.3333
$ws.Range("H126").Value = 1922.6154
$ws.Range("I126").Value = 1527.7142
$ws.Range("J126").Value = 2383.3333
$ws.Range("K126").Value = 4583.142599999999
$ws.Range("L126").Value = 7149.999899999999
$ws.Range("M126").Value = -2113.142599999999
$ws.Range("N126").Value = -12089.9999
$ws.Range("H135").Value = 52020
$ws.Range("J135").Value = 52020
$ws.Range("L135").Value = 52020
$ws.Range("N135").Value = -62160
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27.434782
$ws.Range("I12").Value = 13.818182
$ws.Range("J12").Value = 39.916668
$ws.Range("K12").Value = 41.454546
$ws.Range("L12").Value = 119.750004
$ws.Range("M12").Value = 131.545454
$ws.Range("N12").Value = -465.750004
$ws.Range("H34").Value = 1262.4584
$ws.Range("J34").Value = 1803.0625
$ws.Range("L34").Value = 5409.1875
$ws.Range("N34").Value = -5577.1875
$ws.Range("H46").Value = 1788.8889
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 2128.5715
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 6385.7145
$ws.Range("M46").Value = -1709
$ws.Range("N46").Value = -6567.7145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1752.44
$ws.Range("I97").Value = 1660
$ws.Range("J97").Value = 2122.2
$ws.Range("K97").Value = 1660
$ws.Range("L97").Value = 2122.2
$ws.Range("M97").Value = -1164
$ws.Range("N97").Value = -3114.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5465.2
$ws.Range("I40").Value = 7567.1113
$ws.Range("J40").Value = 3745.4546
$ws.Range("K40").Value = 7567.1113
$ws.Range("L40").Value = 3745.4546
$ws.Range("M40").Value = -7431.1113
$ws.Range("N40").Value = -4017.4546
$ws.Range("H75").Value = 33173
$ws.Range("J75").Value = 33173
$ws.Range("L75").Value = 33173
$ws.Range("N75").Value = -35045
$ws.Range("H78").Value = 33173
$ws.Range("J78").Value = 33173
$ws.Range("L78").Value = 99519
$ws.Range("N78").Value = -108879
$ws.Range("H132").Value = 10008847
$ws.Range("I132").Value = 6262.852
$ws.Range("J132").Value = 21751012
$ws.Range("K132").Value = 18788.556
$ws.Range("L132").Value = 65253036
$ws.Range("M132").Value = -16258.556
$ws.Range("N132").Value = -65258096
$ws.Range("H136").Value = 11114390
$ws.Range("I136").Value = 13514452
$ws.Range("J136").Value = 14100.625
$ws.Range("K136").Value = 40543356
$ws.Range("L136").Value = 42301.875
$ws.Range("M136").Value = -40540806
$ws.Range("N136").Value = -47401.875
$ws.Range("H140").Value = 47744.816
$ws.Range("J140").Value = 47744.816
$ws.Range("L140").Value = 47744.816
$ws.Range("N140").Value = -58104.816
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 47025
$ws.Range("J14").Value = 36050
$ws.Range("L14").Value = 36050
$ws.Range("N14").Value = -36386
$ws.Range("H122").Value = 2177.8462
$ws.Range("I122").Value = 2590.2222
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 7770.6666
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -5320.6666
$ws.Range("N122").Value = -8650
$ws.Range("H126").Value = 2888.318
$ws.Range("I126").Value = 2211.4707
$ws.Range("J126").Value = 5189.6
$ws.Range("K126").Value = 6634.4121
$ws.Range("L126").Value = 15568.8
$ws.Range("M126").Value = -4164.4121
$ws.Range("N126").Value = -20508.8
$ws.Range("H132").Value = 1692.0731
$ws.Range("I132").Value = 1595.2222
$ws.Range("J132").Value = 1878.8572
$ws.Range("K132").Value = 4785.6666
$ws.Range("L132").Value = 5636.571599999999
$ws.Range("M132").Value = -2255.6666
$ws.Range("N132").Value = -10696.5716
$ws.Range("H136").Value = 1139.6471
$ws.Range("I136").Value = 944.1111
$ws.Range("J136").Value = 1893.8572
$ws.Range("K136").Value = 2832.3333
$ws.Range("L136").Value = 5681.571599999999
$ws.Range("M136").Value = -282.3332999999998
$ws.Range("N136").Value = -10781.5716
